$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Round the coordinate values to whole numbers
$ws.Range("Q4").Value = 692986
$ws.Range("R4").Value = 6697797

# Clear the time cells entirely (Starttid / Sluttid)
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
